# Applies the "tr_script synchronized only to tr" edit to the events workbook:
#  - tightens the 5 data-column widths
#  - rewrites the TR-synced Absolute/Relative/Difference timings for TR 1-4 (rows 2-5)
#  - resets rows 6-16 (former TR 5-15) to the blank/zeroed placeholder state,
#    dropping their "Active Stimuli" description

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- column widths (cols A:E) -------------------------------------------------
# The engine's ColumnWidth setter only round-trips to 1/6-character pixel
# granularity, so these are the closest achievable inputs to the target
# stored widths (3.109375 / 12.109375 / 11.5546875 / 11.5546875 / 11.88671875).
$ws.Columns.Item(1).ColumnWidth = 2.3333333333333335
$ws.Columns.Item(2).ColumnWidth = 11.333333333333334
$ws.Columns.Item(3).ColumnWidth = 10.666666666666666
$ws.Columns.Item(4).ColumnWidth = 10.666666666666666
$ws.Columns.Item(5).ColumnWidth = 11.0

# --- rows 2-5 (TR 1-4): refreshed timing values, description unchanged -------
$ws.Range("B2").Value = 1.9998282999986259
$ws.Range("C2").Value = 22520.735716199997
$ws.Range("D2").Value = 22524.737346099999

$ws.Range("B3").Value = 3.9993566000011924
$ws.Range("C3").Value = 22520.7359922
$ws.Range("D3").Value = 22526.736874400001

$ws.Range("B4").Value = 5.9990049000007275
$ws.Range("C4").Value = 22520.7361217
$ws.Range("D4").Value = 22528.736522700001

$ws.Range("B5").Value = 7.9992920000004233
$ws.Range("C5").Value = 22520.736121599999
$ws.Range("D5").Value = 22530.736809800001

# --- rows 6-16 (former TR 5-15): zeroed out, description cleared -------------
for ($r = 6; $r -le 16; $r++) {
    $ws.Range("A$r").Value = 0
    $ws.Range("B$r").Value = 0
    $ws.Range("C$r").Value = 0
    $ws.Range("D$r").Value = 0
    $ws.Range("E$r").ClearContents()
    $ws.Range("E$r").Style = "Normal"
}
